$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 870.53845
$ws.Range("I19").Value = 1033.5
$ws.Range("J19").Value = 730.8570999999999
$ws.Range("K19").Value = 1033.5
$ws.Range("L19").Value = 730.8570999999999
$ws.Range("M19").Value = -858.5
$ws.Range("N19").Value = -1080.8571
# Row 40
$ws.Range("H40").Value = 7389.5
$ws.Range("J40").Value = 9982.5
$ws.Range("L40").Value = 9982.5
$ws.Range("N40").Value = -10332.5
# Row 106
$ws.Range("H106").Value = 4958
$ws.Range("I106").Value = 4958
$ws.Range("K106").Value = 4958
$ws.Range("M106").Value = -4327
# Row 137
$ws.Range("H137").Value = 3000
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -14100

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1455.5385
$ws.Range("I2").Value = 836.8
$ws.Range("J2").Value = 3518
$ws.Range("K2").Value = 836.8
$ws.Range("L2").Value = 3518
$ws.Range("M2").Value = -723.8
$ws.Range("N2").Value = -3744
# Row 74
$ws.Range("H74").Value = 6326.364
$ws.Range("I74").Value = 6421.3335
$ws.Range("K74").Value = 6421.3335
$ws.Range("M74").Value = -5547.3335
# Row 77
$ws.Range("H77").Value = 6326.364
$ws.Range("I77").Value = 6421.3335
$ws.Range("K77").Value = 32106.6675
$ws.Range("M77").Value = -27738.6675
# Row 97
$ws.Range("H97").Value = 616.63635
$ws.Range("I97").Value = 644.8
$ws.Range("K97").Value = 644.8
$ws.Range("M97").Value = -148.8
# Row 102
$ws.Range("H102").Value = 3180.4707
$ws.Range("I102").Value = 1472.3334
$ws.Range("J102").Value = 7280
$ws.Range("K102").Value = 1472.3334
$ws.Range("L102").Value = 7280
$ws.Range("M102").Value = 149.6666
$ws.Range("N102").Value = -10524
# Row 103
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
# Row 110
$ws.Range("H110").Value = 975.3077
$ws.Range("I110").Value = 816.2727
$ws.Range("K110").Value = 816.2727
$ws.Range("M110").Value = 1228.7273
# Row 116
$ws.Range("H116").Value = 1455.5385
$ws.Range("I116").Value = 836.8
$ws.Range("J116").Value = 3518
$ws.Range("K116").Value = 836.8
$ws.Range("L116").Value = 3518
$ws.Range("M116").Value = 1457.2
$ws.Range("N116").Value = -8106
# Row 122
$ws.Range("H122").Value = 1097.8
$ws.Range("I122").Value = 1053.1111
$ws.Range("K122").Value = 3159.3333
$ws.Range("M122").Value = -709.3333000000002
# Row 132
$ws.Range("H132").Value = 1665
$ws.Range("I132").Value = 1665
$ws.Range("K132").Value = 4995
$ws.Range("M132").Value = -2465

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1455.5385
$ws.Range("I3").Value = 836.8
$ws.Range("J3").Value = 3518
$ws.Range("K3").Value = 836.8
$ws.Range("L3").Value = 3518
$ws.Range("M3").Value = -722.8
$ws.Range("N3").Value = -3746
# Row 75
$ws.Range("H75").Value = 16199.75
$ws.Range("I75").Value = 3266.3333
$ws.Range("J75").Value = 55000
$ws.Range("K75").Value = 3266.3333
$ws.Range("L75").Value = 55000
$ws.Range("M75").Value = -2330.3333
$ws.Range("N75").Value = -56872
# Row 78
$ws.Range("H78").Value = 16199.75
$ws.Range("I78").Value = 3266.3333
$ws.Range("J78").Value = 55000
$ws.Range("K78").Value = 9798.999899999999
$ws.Range("L78").Value = 165000
$ws.Range("M78").Value = -5118.999899999999
$ws.Range("N78").Value = -174360
# Row 86
$ws.Range("H86").Value = 3875.85
$ws.Range("I86").Value = 1393.4166
$ws.Range("K86").Value = 1393.4166
$ws.Range("M86").Value = -270.4166
# Row 89
$ws.Range("H89").Value = 3875.85
$ws.Range("I89").Value = 1393.4166
$ws.Range("K89").Value = 6967.083000000001
$ws.Range("M89").Value = -1351.083000000001
# Row 107
$ws.Range("H107").Value = 6236.364
$ws.Range("I107").Value = 4825
$ws.Range("K107").Value = 4825
$ws.Range("M107").Value = -2905
# Row 134
$ws.Range("H134").Value = 1000
$ws.Range("I134").Value = 1000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -465
$ws.Range("N134").Value = ""

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 43
$ws.Range("H43").Value = 20409.285
$ws.Range("J43").Value = 20409.285
$ws.Range("L43").Value = 20409.285
$ws.Range("N43").Value = -20777.285
# Row 58
$ws.Range("H58").Value = 5135.2
$ws.Range("I58").Value = 4629.5713
$ws.Range("K58").Value = 4629.5713
$ws.Range("M58").Value = -4426.5713
# Row 74
$ws.Range("H74").Value = 45225.6
$ws.Range("J74").Value = 46917.332
$ws.Range("L74").Value = 46917.332
$ws.Range("N74").Value = -48665.332
# Row 77
$ws.Range("H77").Value = 45225.6
$ws.Range("J77").Value = 46917.332
$ws.Range("L77").Value = 140751.996
$ws.Range("N77").Value = -149487.996
# Row 101
$ws.Range("H101").Value = 20409.285
$ws.Range("J101").Value = 20409.285
$ws.Range("L101").Value = 20409.285
$ws.Range("N101").Value = -26899.285
# Row 107
$ws.Range("H107").Value = 415.36365
$ws.Range("I107").Value = 287.5
$ws.Range("J107").Value = 568.8
$ws.Range("K107").Value = 287.5
$ws.Range("L107").Value = 568.8
$ws.Range("M107").Value = 1632.5
$ws.Range("N107").Value = -4408.8
# Row 122
$ws.Range("H122").Value = 1912
$ws.Range("I122").Value = 1912
$ws.Range("K122").Value = 5736
$ws.Range("M122").Value = -3286
# Row 136
$ws.Range("H136").Value = 5135.2
$ws.Range("I136").Value = 4629.5713
$ws.Range("K136").Value = 13888.7139
$ws.Range("M136").Value = -11338.7139

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 103.23529
$ws.Range("J12").Value = 124.76923
$ws.Range("L12").Value = 374.30769
$ws.Range("N12").Value = -720.30769
# Row 14
$ws.Range("H14").Value = 543.4286
$ws.Range("I14").Value = 543.4286
$ws.Range("K14").Value = 1630.2858
$ws.Range("M14").Value = -1457.2858

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 292.4091
$ws.Range("I2").Value = 152.07692
$ws.Range("K2").Value = 152.07692
$ws.Range("M2").Value = -39.07692
# Row 36
$ws.Range("H36").Value = 1831.1666
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4515
# Row 42
$ws.Range("H42").Value = 99999
$ws.Range("J42").Value = 99999
$ws.Range("L42").Value = 99999
$ws.Range("N42").Value = -100969
# Row 97
$ws.Range("H97").Value = 680.8333
$ws.Range("I97").Value = 657
$ws.Range("K97").Value = 657
$ws.Range("M97").Value = -161
# Row 104
$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -26988
# Row 107
$ws.Range("H107").Value = 2486.2222
$ws.Range("I107").Value = 4091
$ws.Range("J107").Value = 1683.8334
$ws.Range("K107").Value = 4091
$ws.Range("L107").Value = 1683.8334
$ws.Range("M107").Value = -2171
$ws.Range("N107").Value = -5523.8334
# Row 115
$ws.Range("H115").Value = 99999
$ws.Range("J115").Value = 99999
$ws.Range("L115").Value = 99999
$ws.Range("N115").Value = -102349

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 6843.5
$ws.Range("I68").Value = 3812
$ws.Range("K68").Value = 3812
$ws.Range("M68").Value = -3063
# Row 71
$ws.Range("H71").Value = 6843.5
$ws.Range("I71").Value = 3812
$ws.Range("K71").Value = 19060
$ws.Range("M71").Value = -15316
# Row 136
$ws.Range("H136").Value = 3498.5
$ws.Range("J136").Value = 3498.5
$ws.Range("L136").Value = 10495.5
$ws.Range("N136").Value = -15595.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 23500
$ws.Range("J54").Value = 23500
$ws.Range("L54").Value = 23500
$ws.Range("N54").Value = -24540
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").Value = ""
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").Value = ""
# Row 92
$ws.Range("H92").Value = 16678
$ws.Range("J92").Value = 16678
$ws.Range("L92").Value = 16678
$ws.Range("N92").Value = -21670
# Row 113
$ws.Range("H113").Value = 811.2
$ws.Range("I113").Value = 701.625
$ws.Range("K113").Value = 2104.875
$ws.Range("M113").Value = 65.125
# Row 136
$ws.Range("H136").Value = 4156.25
$ws.Range("I136").Value = 2303
$ws.Range("J136").Value = 6346.4546
$ws.Range("K136").Value = 6909
$ws.Range("L136").Value = 19039.3638
$ws.Range("M136").Value = -4359
$ws.Range("N136").Value = -24139.3638
